$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AMSIN")
$ws.Cells.Item(43, 1).Value = "test"
